$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update period (E column) values for rows 16-73 to reflect reordering of periods
$ws.Range("E16").Value = "2103"
$ws.Range("E17").Value = "2102"
$ws.Range("E18").Value = "2101"
$ws.Range("E19").Value = "2012"
$ws.Range("E20").Value = "2011"
$ws.Range("E21").Value = "2010"
$ws.Range("E22").Value = "2009"
$ws.Range("E23").Value = "2008"
$ws.Range("E24").Value = "2007"
$ws.Range("E25").Value = "2006"
$ws.Range("E26").Value = "2005"
$ws.Range("E27").Value = "2004"
$ws.Range("E28").Value = "2003"
$ws.Range("E29").Value = "2002"
$ws.Range("E30").Value = "2001"
$ws.Range("E31").Value = "1912"
$ws.Range("E32").Value = "1911"
$ws.Range("E33").Value = "1910"
$ws.Range("E34").Value = "1909"
$ws.Range("E35").Value = "1908"
$ws.Range("E36").Value = "1907"
$ws.Range("E37").Value = "1906"
$ws.Range("E38").Value = "1905"
$ws.Range("E39").Value = "1904"
$ws.Range("E40").Value = "1903"
$ws.Range("E41").Value = "1902"
$ws.Range("E42").Value = "1901"
$ws.Range("E43").Value = "1812"
$ws.Range("E44").Value = "1811"
$ws.Range("E45").Value = "1810"
$ws.Range("E46").Value = "1809"
$ws.Range("E47").Value = "1808"
$ws.Range("E48").Value = "1807"
$ws.Range("E49").Value = "1806"
$ws.Range("E50").Value = "1805"
$ws.Range("E51").Value = "1804"
$ws.Range("E52").Value = "1803"
$ws.Range("E53").Value = "1802"
$ws.Range("E54").Value = "1801"
$ws.Range("E55").Value = "1712"
$ws.Range("E56").Value = "1711"
$ws.Range("E57").Value = "1710"
$ws.Range("E58").Value = "1709"
$ws.Range("E59").Value = "1708"
$ws.Range("E60").Value = "1707"
$ws.Range("E61").Value = "1706"
$ws.Range("E62").Value = "1705"
$ws.Range("E63").Value = "1704"
$ws.Range("E64").Value = "1703"
$ws.Range("E65").Value = "1702"
$ws.Range("E66").Value = "1701"
$ws.Range("E67").Value = "1612"
$ws.Range("E68").Value = "1611"
$ws.Range("E69").Value = "1610"
$ws.Range("E70").Value = "1609"
$ws.Range("E71").Value = "1608"
$ws.Range("E72").Value = "1607"
$ws.Range("E73").Value = "1711"

# Row 16: Valor Mora (F) changes
$ws.Range("F16").Value = 39866

# Row 33 becomes a normal LEDIS MABEL row (previously held CARLOS IVAN data)
$ws.Range("C33").Value = "60315332"
$ws.Range("D33").Value = "LEDIS MABEL MARIA VERGARA ALVAREZ"
$ws.Range("F33").Value = 52000
$ws.Range("G33").Value = 1300000

# Row 73 now holds CARLOS IVAN RUA SERRANO data (moved from row 33)
$ws.Range("C73").Value = "1143374662"
$ws.Range("D73").Value = "CARLOS IVAN RUA SERRANO"
$ws.Range("F73").Value = 36000
$ws.Range("G73").Value = 900000
